$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update defined names: add hidden _FilterDatabase name, extend DOLAR_OBS_ADO range ---
$ws.Names.Add("_xlnm._FilterDatabase", "=DOLAR_OBS_ADO!`$A`$3:`$B`$762") | Out-Null
$ws.Names.Item(1).Visible = $false
$wb.Names.Item(2).RefersTo = "=DOLAR_OBS_ADO!`$A`$1:`$B`$762"

# --- Bulk-apply existing row formatting (numeric style) to the whole new range 684:762 ---
$ws.Range("A683:B683").Copy() | Out-Null
$ws.Range("A684:B762").PasteSpecial(-4122) | Out-Null

# --- Write date + value data for the new rows ---
$ws.Cells.Item(684, 1).Value = 44147
$ws.Cells.Item(684, 2).Value = 757.42
$ws.Cells.Item(685, 1).Value = 44148
$ws.Cells.Item(685, 2).Value = 757.43
$ws.Cells.Item(686, 1).Value = 44149
$ws.Cells.Item(686, 2).Value = "--"
$ws.Cells.Item(687, 1).Value = 44150
$ws.Cells.Item(687, 2).Value = "--"
$ws.Cells.Item(688, 1).Value = 44151
$ws.Cells.Item(688, 2).Value = 766.7
$ws.Cells.Item(689, 1).Value = 44152
$ws.Cells.Item(689, 2).Value = 767.86
$ws.Cells.Item(690, 1).Value = 44153
$ws.Cells.Item(690, 2).Value = 767.05
$ws.Cells.Item(691, 1).Value = 44154
$ws.Cells.Item(691, 2).Value = 758.1
$ws.Cells.Item(692, 1).Value = 44155
$ws.Cells.Item(692, 2).Value = 758.62
$ws.Cells.Item(693, 1).Value = 44156
$ws.Cells.Item(693, 2).Value = "--"
$ws.Cells.Item(694, 1).Value = 44157
$ws.Cells.Item(694, 2).Value = "--"
$ws.Cells.Item(695, 1).Value = 44158
$ws.Cells.Item(695, 2).Value = 761.55
$ws.Cells.Item(696, 1).Value = 44159
$ws.Cells.Item(696, 2).Value = 765.96
$ws.Cells.Item(697, 1).Value = 44160
$ws.Cells.Item(697, 2).Value = 772.83
$ws.Cells.Item(698, 1).Value = 44161
$ws.Cells.Item(698, 2).Value = 771.68
$ws.Cells.Item(699, 1).Value = 44162
$ws.Cells.Item(699, 2).Value = 766
$ws.Cells.Item(700, 1).Value = 44163
$ws.Cells.Item(700, 2).Value = "--"
$ws.Cells.Item(701, 1).Value = 44164
$ws.Cells.Item(701, 2).Value = "--"
$ws.Cells.Item(702, 1).Value = 44165
$ws.Cells.Item(702, 2).Value = 766.69
$ws.Cells.Item(703, 1).Value = 44166
$ws.Cells.Item(703, 2).Value = 767.29
$ws.Cells.Item(704, 1).Value = 44167
$ws.Cells.Item(704, 2).Value = 760.16
$ws.Cells.Item(705, 1).Value = 44168
$ws.Cells.Item(705, 2).Value = 755.34
$ws.Cells.Item(706, 1).Value = 44169
$ws.Cells.Item(706, 2).Value = 752.03
$ws.Cells.Item(707, 1).Value = 44170
$ws.Cells.Item(707, 2).Value = "--"
$ws.Cells.Item(708, 1).Value = 44171
$ws.Cells.Item(708, 2).Value = "--"
$ws.Cells.Item(709, 1).Value = 44172
$ws.Cells.Item(709, 2).Value = 747.61
$ws.Cells.Item(710, 1).Value = 44173
$ws.Cells.Item(710, 2).Value = "--"
$ws.Cells.Item(711, 1).Value = 44174
$ws.Cells.Item(711, 2).Value = 744.82
$ws.Cells.Item(712, 1).Value = 44175
$ws.Cells.Item(712, 2).Value = 739.45
$ws.Cells.Item(713, 1).Value = 44176
$ws.Cells.Item(713, 2).Value = 738.17
$ws.Cells.Item(714, 1).Value = 44177
$ws.Cells.Item(714, 2).Value = "--"
$ws.Cells.Item(715, 1).Value = 44178
$ws.Cells.Item(715, 2).Value = "--"
$ws.Cells.Item(716, 1).Value = 44179
$ws.Cells.Item(716, 2).Value = 733.55
$ws.Cells.Item(717, 1).Value = 44180
$ws.Cells.Item(717, 2).Value = 731.58
$ws.Cells.Item(718, 1).Value = 44181
$ws.Cells.Item(718, 2).Value = 734.23
$ws.Cells.Item(719, 1).Value = 44182
$ws.Cells.Item(719, 2).Value = 735.09
$ws.Cells.Item(720, 1).Value = 44183
$ws.Cells.Item(720, 2).Value = 723.44
$ws.Cells.Item(721, 1).Value = 44184
$ws.Cells.Item(721, 2).Value = "--"
$ws.Cells.Item(722, 1).Value = 44185
$ws.Cells.Item(722, 2).Value = "--"
$ws.Cells.Item(723, 1).Value = 44186
$ws.Cells.Item(723, 2).Value = 723.85
$ws.Cells.Item(724, 1).Value = 44187
$ws.Cells.Item(724, 2).Value = 730.7
$ws.Cells.Item(725, 1).Value = 44188
$ws.Cells.Item(725, 2).Value = 728.96
$ws.Cells.Item(726, 1).Value = 44189
$ws.Cells.Item(726, 2).Value = 716.25
$ws.Cells.Item(727, 1).Value = 44190
$ws.Cells.Item(727, 2).Value = "--"
$ws.Cells.Item(728, 1).Value = 44191
$ws.Cells.Item(728, 2).Value = "--"
$ws.Cells.Item(729, 1).Value = 44192
$ws.Cells.Item(729, 2).Value = "--"
$ws.Cells.Item(730, 1).Value = 44193
$ws.Cells.Item(730, 2).Value = 710.26
$ws.Cells.Item(731, 1).Value = 44194
$ws.Cells.Item(731, 2).Value = 710.64
$ws.Cells.Item(732, 1).Value = 44195
$ws.Cells.Item(732, 2).Value = 711.24
$ws.Cells.Item(733, 1).Value = 44196
$ws.Cells.Item(733, 2).Value = "--"
$ws.Cells.Item(734, 1).Value = 44197
$ws.Cells.Item(734, 2).Value = "--"
$ws.Cells.Item(735, 1).Value = 44198
$ws.Cells.Item(735, 2).Value = "--"
$ws.Cells.Item(736, 1).Value = 44199
$ws.Cells.Item(736, 2).Value = "--"
$ws.Cells.Item(737, 1).Value = 44200
$ws.Cells.Item(737, 2).Value = 710.95
$ws.Cells.Item(738, 1).Value = 44201
$ws.Cells.Item(738, 2).Value = 702.93
$ws.Cells.Item(739, 1).Value = 44202
$ws.Cells.Item(739, 2).Value = 702.29
$ws.Cells.Item(740, 1).Value = 44203
$ws.Cells.Item(740, 2).Value = 696.18
$ws.Cells.Item(741, 1).Value = 44204
$ws.Cells.Item(741, 2).Value = 709.99
$ws.Cells.Item(742, 1).Value = 44205
$ws.Cells.Item(742, 2).Value = "--"
$ws.Cells.Item(743, 1).Value = 44206
$ws.Cells.Item(743, 2).Value = "--"
$ws.Cells.Item(744, 1).Value = 44207
$ws.Cells.Item(744, 2).Value = 713.28
$ws.Cells.Item(745, 1).Value = 44208
$ws.Cells.Item(745, 2).Value = 718.89
$ws.Cells.Item(746, 1).Value = 44209
$ws.Cells.Item(746, 2).Value = 725.24
$ws.Cells.Item(747, 1).Value = 44210
$ws.Cells.Item(747, 2).Value = 739.72
$ws.Cells.Item(748, 1).Value = 44211
$ws.Cells.Item(748, 2).Value = 735.35
$ws.Cells.Item(749, 1).Value = 44212
$ws.Cells.Item(749, 2).Value = "--"
$ws.Cells.Item(750, 1).Value = 44213
$ws.Cells.Item(750, 2).Value = "--"
$ws.Cells.Item(751, 1).Value = 44214
$ws.Cells.Item(751, 2).Value = 735.06
$ws.Cells.Item(752, 1).Value = 44215
$ws.Cells.Item(752, 2).Value = 736.11
$ws.Cells.Item(753, 1).Value = 44216
$ws.Cells.Item(753, 2).Value = 733.73
$ws.Cells.Item(754, 1).Value = 44217
$ws.Cells.Item(754, 2).Value = 730.38
$ws.Cells.Item(755, 1).Value = 44218
$ws.Cells.Item(755, 2).Value = 715.56
$ws.Cells.Item(756, 1).Value = 44219
$ws.Cells.Item(756, 2).Value = "--"
$ws.Cells.Item(757, 1).Value = 44220
$ws.Cells.Item(757, 2).Value = "--"
$ws.Cells.Item(758, 1).Value = 44221
$ws.Cells.Item(758, 2).Value = 724.26
$ws.Cells.Item(759, 1).Value = 44222
$ws.Cells.Item(759, 2).Value = 731.92
$ws.Cells.Item(760, 1).Value = 44223
$ws.Cells.Item(760, 2).Value = 731
$ws.Cells.Item(761, 1).Value = 44224
$ws.Cells.Item(761, 2).Value = 736.88
$ws.Cells.Item(762, 1).Value = 44225
$ws.Cells.Item(762, 2).Value = 741.4

# --- Fix style for dash ("--") cells to match existing dash-cell style (copy from B7) ---
$ws.Range("A7:B7").Copy() | Out-Null
$ws.Range("A686:B686").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(686, 1).Value = 44149
$ws.Cells.Item(686, 2).Value = "--"
$ws.Range("A687:B687").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(687, 1).Value = 44150
$ws.Cells.Item(687, 2).Value = "--"
$ws.Range("A693:B693").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(693, 1).Value = 44156
$ws.Cells.Item(693, 2).Value = "--"
$ws.Range("A694:B694").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(694, 1).Value = 44157
$ws.Cells.Item(694, 2).Value = "--"
$ws.Range("A700:B700").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(700, 1).Value = 44163
$ws.Cells.Item(700, 2).Value = "--"
$ws.Range("A701:B701").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(701, 1).Value = 44164
$ws.Cells.Item(701, 2).Value = "--"
$ws.Range("A707:B707").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(707, 1).Value = 44170
$ws.Cells.Item(707, 2).Value = "--"
$ws.Range("A708:B708").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(708, 1).Value = 44171
$ws.Cells.Item(708, 2).Value = "--"
$ws.Range("A710:B710").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(710, 1).Value = 44173
$ws.Cells.Item(710, 2).Value = "--"
$ws.Range("A714:B714").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(714, 1).Value = 44177
$ws.Cells.Item(714, 2).Value = "--"
$ws.Range("A715:B715").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(715, 1).Value = 44178
$ws.Cells.Item(715, 2).Value = "--"
$ws.Range("A721:B721").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(721, 1).Value = 44184
$ws.Cells.Item(721, 2).Value = "--"
$ws.Range("A722:B722").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(722, 1).Value = 44185
$ws.Cells.Item(722, 2).Value = "--"
$ws.Range("A727:B727").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(727, 1).Value = 44190
$ws.Cells.Item(727, 2).Value = "--"
$ws.Range("A728:B728").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(728, 1).Value = 44191
$ws.Cells.Item(728, 2).Value = "--"
$ws.Range("A729:B729").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(729, 1).Value = 44192
$ws.Cells.Item(729, 2).Value = "--"
$ws.Range("A733:B733").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(733, 1).Value = 44196
$ws.Cells.Item(733, 2).Value = "--"
$ws.Range("A734:B734").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(734, 1).Value = 44197
$ws.Cells.Item(734, 2).Value = "--"
$ws.Range("A735:B735").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(735, 1).Value = 44198
$ws.Cells.Item(735, 2).Value = "--"
$ws.Range("A736:B736").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(736, 1).Value = 44199
$ws.Cells.Item(736, 2).Value = "--"
$ws.Range("A742:B742").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(742, 1).Value = 44205
$ws.Cells.Item(742, 2).Value = "--"
$ws.Range("A743:B743").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(743, 1).Value = 44206
$ws.Cells.Item(743, 2).Value = "--"
$ws.Range("A749:B749").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(749, 1).Value = 44212
$ws.Cells.Item(749, 2).Value = "--"
$ws.Range("A750:B750").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(750, 1).Value = 44213
$ws.Cells.Item(750, 2).Value = "--"
$ws.Range("A756:B756").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(756, 1).Value = 44219
$ws.Cells.Item(756, 2).Value = "--"
$ws.Range("A757:B757").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(757, 1).Value = 44220
$ws.Cells.Item(757, 2).Value = "--"

# --- Update column widths to best match target ---
$ws.Columns.Item(1).ColumnWidth = 11.877604166666666
$ws.Columns.Item(2).ColumnWidth = 13.592447916666666